$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.549.04"
$ws.Range("E2").Value = "  -3.63%  "
$ws.Range("D3").Value = "2.999.04"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'537.33"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'132.86"
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "2.992.80"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  -5.08%  "
$ws.Range("D11").Value = "'6.09"
$ws.Range("E11").Value = "  -5.58%  "
$ws.Range("D12").Value = "'0.447"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "'0.0000222"
$ws.Range("E13").Value = "  -2.06%  "
$ws.Range("D14").Value = "'33.84"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "3.484.26"
$ws.Range("E15").Value = "  -2.85%  "
$ws.Range("D16").Value = "61.615.23"
$ws.Range("E16").Value = "  -3.66%  "
$ws.Range("D17").Value = "'0.110"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "3.002.33"
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("D19").Value = "'6.60"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").Value = "'470.77"
$ws.Range("E20").Value = "  -2.58%  "
$ws.Range("D21").Value = "'13.22"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "'0.671"
$ws.Range("E22").Value = "  -4.54%  "
$ws.Range("D23").Value = "'6.93"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").Value = "'80.28"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").Value = "'11.97"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'7.77"
$ws.Range("E28").Value = "  -4.49%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.88"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("B31").Value = "Mantle"
$ws.Range("C31").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'25.54"
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.30"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'5.46"
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'55.32"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("D36").Value = "'5.89"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").Value = "'461.16"
$ws.Range("E37").Value = "  -8.68%  "
$ws.Range("D38").Value = "3.175.76"
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("D39").Value = "'0.0789"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.120"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0384"
$ws.Range("E41").Value = "  -3.60%  "
$ws.Range("D42").Value = "'8.08"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'2.43"
$ws.Range("E43").Value = "  -7.44%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "'25.81"
$ws.Range("E45").Value = "  +5.05%  "
$ws.Range("D46").Value = "'0.244"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("D47").Value = "'2.00"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").Value = "'117.19"
$ws.Range("E49").Value = "  -3.96%  "
$ws.Range("D50").Value = "0.0₃0492"
$ws.Range("E50").Value = "  -7.51%  "
$ws.Range("E51").Value = "  +5.93%  "
